$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume refresh - GitHub Actions scheduled update

$ws.Range("D2").Value = "'30.652.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").Value = "'2.114.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("E4").Value = "  +1.01%  "

$ws.Range("D5").Value = "'349.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.85%  "

$ws.Range("D6").Value = "'1.010"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.87%  "

$ws.Range("D7").Value = "'0.5261"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.35%  "

$ws.Range("D8").Value = "'0.4511"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.63%  "

$ws.Range("D9").Value = "'53.59"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "

$ws.Range("D10").Value = "'0.09033"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.95%  "

$ws.Range("E11").Value = "  -0.50%  "

$ws.Range("D12").Value = "'24.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.22%  "

$ws.Range("D13").Value = "'2.113.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.16%  "

$ws.Range("D14").Value = "'6.817"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.25%  "

$ws.Range("D15").Value = "'8.036"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.89%  "

$ws.Range("D16").Value = "'100.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.02%  "

$ws.Range("D17").Value = "'0.00001170"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.27%  "

$ws.Range("D18").Value = "'1.012"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.94%  "

$ws.Range("E19").Value = "  +1.08%  "

$ws.Range("D20").Value = "'19.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("D21").Value = "'1.010"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.97%  "

$ws.Range("D22").Value = "'6.295"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.15%  "

$ws.Range("D23").Value = "'30.719.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.66%  "

$ws.Range("D24").Value = "'12.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.99%  "

$ws.Range("D25").Value = "'2.388"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.05%  "

$ws.Range("D26").Value = "'2.354.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.80%  "

$ws.Range("D27").Value = "'22.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.20%  "

$ws.Range("D28").Value = "'165.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.28%  "

$ws.Range("D29").Value = "'2.533"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.43%  "

$ws.Range("D30").Value = "'135.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.11%  "

$ws.Range("D31").Value = "'1.190"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.95%  "

$ws.Range("D32").Value = "'0.1074"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("D33").Value = "'1.651"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.34%  "

$ws.Range("D34").Value = "'6.362"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.50%  "

$ws.Range("D35").Value = "'4.013"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.25%  "

$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'5.922"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.59%  "

$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "'10.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.28%  "

$ws.Range("D38").Value = "'0.02651"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.83%  "

$ws.Range("D39").Value = "'0.06842"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.05%  "

$ws.Range("D40").Value = "'0.2313"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.34%  "

$ws.Range("D41").Value = "'12.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.14%  "

$ws.Range("E42").Value = "  -0.21%  "

$ws.Range("D43").Value = "'1.280"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.78%  "

$ws.Range("D44").Value = "'14.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.31%  "

$ws.Range("D45").Value = "'0.6439"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.80%  "

$ws.Range("D46").Value = "'2.323"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.32%  "

$ws.Range("D47").Value = "'3.753"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.36%  "

$ws.Range("D48").Value = "'0.00000000354"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.62%  "

$ws.Range("D49").Value = "'1.250"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.12%  "

$ws.Range("D50").Value = "'0.07279"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.15%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'82.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.54%  "

Write-Host "Applied cryptos update"
